$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 6-11: the nine data-bearing columns (A,B,D,E,F,G,H,Q,R) got
# cyclically re-permuted among themselves (rows 7,9,10,11 rotate; rows
# 6 and 8 swap). All other columns in these rows stay as-is.

function Set-RowData {
    param($row, $A, $B, $D, $E, $F, $G, $H, $Q, $R)
    $ws.Cells.Item($row, 1).Value  = $A
    $ws.Cells.Item($row, 2).Value  = $B
    $ws.Cells.Item($row, 4).Value  = $D
    $ws.Cells.Item($row, 5).Value  = $E
    $ws.Cells.Item($row, 6).Value  = $F
    $ws.Cells.Item($row, 7).Value  = $G
    $ws.Cells.Item($row, 8).Value  = $H
    $ws.Cells.Item($row, 17).Value = $Q
    $ws.Cells.Item($row, 18).Value = $R
}

Set-RowData 6  130803069 81228 "NT" 1049   "Kortskaftad ärgspik" "Microcalicium ahlneri"  "Tibell"             424815 6712165
Set-RowData 7  130803088 83089 "NT" 1312   "Gammelgransskål"     "Pseudographis pinicola" "(Nyl.) Rehm"        424964 6712067
Set-RowData 8  130803067 78255 "NT" 228579 "Liten svartspik"     "Chaenothecopsis nana"   "Tibell"             424814 6712361
Set-RowData 9  130803042 91771 "LC" 5447   "Vedticka"            "Fuscoporia viticola"    "(Schwein.) Murrill" 424979 6712092
Set-RowData 10 130803071 91181 "LC" 5685   "Gullgröppa"          "Pseudomerulius aureus"  "(Fr.) Jülich"       424873 6712126
Set-RowData 11 130803064 91829 "NT" 5442   "Tallticka"           "Porodaedalea pini"      "(Brot.) Murrill"    424893 6712101

# New row 74, appended after the previous last row (73).
$r = 74
$ws.Cells.Item($r, 1).Value  = 130849489
$ws.Cells.Item($r, 2).Value  = 80349
$ws.Cells.Item($r, 4).Value  = "NT"
$ws.Cells.Item($r, 5).Value  = 2081
$ws.Cells.Item($r, 6).Value  = "Skrovellav"
$ws.Cells.Item($r, 7).Value  = "Lobaria scrobiculata"
$ws.Cells.Item($r, 8).Value  = "(Scop.) DC."
$ws.Cells.Item($r, 16).Value = "Nordvallen, Dlr"
$ws.Cells.Item($r, 17).Value = 424637
$ws.Cells.Item($r, 18).Value = 6712469
$ws.Cells.Item($r, 19).Value = 20
$ws.Cells.Item($r, 20).Value = "Dalarna"
$ws.Cells.Item($r, 21).Value = "Malung-Sälen"
$ws.Cells.Item($r, 22).Value = "Dalarna"
$ws.Cells.Item($r, 23).Value = "Malung"

# Y/AA hold plain date-text (not real dates) throughout this sheet; force
# text so Excel's auto date-detection doesn't coerce them to serials, then
# drop back to the default style so no stray NumberFormat override sticks.
$ws.Cells.Item($r, 25).NumberFormat = "@"
$ws.Cells.Item($r, 25).Value = "2026-01-23"
$ws.Cells.Item($r, 25).Style = "Normal"
$ws.Cells.Item($r, 27).NumberFormat = "@"
$ws.Cells.Item($r, 27).Value = "2026-01-23"
$ws.Cells.Item($r, 27).Style = "Normal"

$ws.Cells.Item($r, 30).Value = $false
$ws.Cells.Item($r, 31).Value = $false
$ws.Cells.Item($r, 33).Value = $false
$ws.Cells.Item($r, 49).Value = "Moa Björnberg dillner"
$ws.Cells.Item($r, 50).Value = "Moa Björnberg dillner"
